# Refresh market-price derived columns (H:N) across the Leve profit sheets.
# Values come from a scheduled market-data refresh; row identity (G = Leve Item ID) is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4999
$ws.Range("I76").Value = 4999
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4999
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4684
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 4999
$ws.Range("I79").Value = 4999
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4999
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3907
$ws.Range("N79").ClearContents()
$ws.Range("H129").Value = 2415.476
$ws.Range("I129").Value = 2200.4707
$ws.Range("K129").Value = 6601.4121
$ws.Range("M129").Value = -1601.4121
$ws.Range("H137").Value = 2736.9688
$ws.Range("I137").Value = 1809.75
$ws.Range("J137").Value = 3293.3
$ws.Range("K137").Value = 5429.25
$ws.Range("L137").Value = 9879.900000000001
$ws.Range("M137").Value = -2879.25
$ws.Range("N137").Value = -14979.9
$ws.Range("H138").Value = 3524.3103
$ws.Range("J138").Value = 4075.5405
$ws.Range("L138").Value = 12226.6215
$ws.Range("N138").Value = -22506.6215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5953.47
$ws.Range("I32").Value = 3507.0361
$ws.Range("J32").Value = 17897.824
$ws.Range("K32").Value = 3507.0361
$ws.Range("L32").Value = 17897.824
$ws.Range("M32").Value = -3220.0361
$ws.Range("N32").Value = -18471.824
$ws.Range("H45").Value = 3681.2942
$ws.Range("I45").Value = 3598.7856
$ws.Range("K45").Value = 3598.7856
$ws.Range("M45").Value = -3221.7856
$ws.Range("H61").Value = 8773.700000000001
$ws.Range("I61").Value = 3872.9
$ws.Range("J61").Value = 13674.5
$ws.Range("K61").Value = 3872.9
$ws.Range("L61").Value = 13674.5
$ws.Range("M61").Value = -3660.9
$ws.Range("N61").Value = -14098.5
$ws.Range("H74").Value = 11370491
$ws.Range("I74").Value = 41669860
$ws.Range("K74").Value = 41669860
$ws.Range("M74").Value = -41668986
$ws.Range("H77").Value = 11370491
$ws.Range("I77").Value = 41669860
$ws.Range("K77").Value = 208349300
$ws.Range("M77").Value = -208344932
$ws.Range("H132").Value = 780968.9399999999
$ws.Range("I132").Value = 1167258.9
$ws.Range("J132").Value = 8389.1
$ws.Range("K132").Value = 3501776.7
$ws.Range("L132").Value = 25167.3
$ws.Range("M132").Value = -3499246.7
$ws.Range("N132").Value = -30227.3
$ws.Range("H136").Value = 8773.700000000001
$ws.Range("I136").Value = 3872.9
$ws.Range("J136").Value = 13674.5
$ws.Range("K136").Value = 11618.7
$ws.Range("L136").Value = 41023.5
$ws.Range("M136").Value = -9068.700000000001
$ws.Range("N136").Value = -46123.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 30000
$ws.Range("J27").Value = 30000
$ws.Range("L27").Value = 30000
$ws.Range("N27").Value = -30384
$ws.Range("H94").Value = 834.03125
$ws.Range("I94").Value = 667.8570999999999
$ws.Range("K94").Value = 667.8570999999999
$ws.Range("M94").Value = -216.8570999999999
$ws.Range("H134").Value = 714257.3
$ws.Range("I134").Value = 1082318
$ws.Range("J134").Value = 8807.75
$ws.Range("K134").Value = 3246954
$ws.Range("L134").Value = 26423.25
$ws.Range("M134").Value = -3244419
$ws.Range("N134").Value = -31493.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 169.8125
$ws.Range("I7").Value = 88.545456
$ws.Range("J7").Value = 348.6
$ws.Range("K7").Value = 88.545456
$ws.Range("L7").Value = 348.6
$ws.Range("M7").Value = 24.454544
$ws.Range("N7").Value = -574.6
$ws.Range("H50").Value = 119994.75
$ws.Range("J50").Value = 119994.75
$ws.Range("L50").Value = 119994.75
$ws.Range("N50").Value = -121244.75
$ws.Range("H58").Value = 1032333.7
$ws.Range("J58").Value = 9665
$ws.Range("L58").Value = 9665
$ws.Range("N58").Value = -10071
$ws.Range("H68").Value = 66947.2
$ws.Range("J68").Value = 66947.2
$ws.Range("L68").Value = 66947.2
$ws.Range("N68").Value = -68445.2
$ws.Range("H71").Value = 66947.2
$ws.Range("J71").Value = 66947.2
$ws.Range("L71").Value = 200841.6
$ws.Range("N71").Value = -208329.6
$ws.Range("H99").Value = 4932.1665
$ws.Range("J99").Value = 5547.1
$ws.Range("L99").Value = 5547.1
$ws.Range("N99").Value = -8543.1
$ws.Range("H126").Value = 4932.1665
$ws.Range("J126").Value = 5547.1
$ws.Range("L126").Value = 16641.3
$ws.Range("N126").Value = -21581.3
$ws.Range("H136").Value = 1032333.7
$ws.Range("J136").Value = 9665
$ws.Range("L136").Value = 28995
$ws.Range("N136").Value = -34095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2540.6667
$ws.Range("I69").Value = 2425
$ws.Range("J69").Value = 2598.5
$ws.Range("K69").Value = 7275
$ws.Range("L69").Value = 7795.5
$ws.Range("M69").Value = -6464
$ws.Range("N69").Value = -9417.5
$ws.Range("H72").Value = 2540.6667
$ws.Range("I72").Value = 2425
$ws.Range("J72").Value = 2598.5
$ws.Range("K72").Value = 21825
$ws.Range("L72").Value = 23386.5
$ws.Range("M72").Value = -17769
$ws.Range("N72").Value = -31498.5
$ws.Range("H109").Value = 4667.7
$ws.Range("I109").Value = 2327.8572
$ws.Range("K109").Value = 6983.571599999999
$ws.Range("M109").Value = -5943.571599999999
$ws.Range("H122").Value = 767.5
$ws.Range("J122").Value = 831.19354
$ws.Range("L122").Value = 7480.74186
$ws.Range("N122").Value = -12380.74186
$ws.Range("H137").Value = 1938.1818
$ws.Range("I137").Value = 2220.182
$ws.Range("J137").Value = 1656.1818
$ws.Range("K137").Value = 6660.545999999999
$ws.Range("L137").Value = 4968.5454
$ws.Range("M137").Value = -1560.545999999999
$ws.Range("N137").Value = -15168.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1107.875
$ws.Range("I132").Value = 1191.6666
$ws.Range("K132").Value = 3574.9998
$ws.Range("M132").Value = -1044.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4571.7393
$ws.Range("I7").Value = 4620.8887
$ws.Range("K7").Value = 4620.8887
$ws.Range("M7").Value = -4508.8887
$ws.Range("H56").Value = 35000
$ws.Range("J56").Value = 35000
$ws.Range("L56").Value = 35000
$ws.Range("N56").Value = -36382
$ws.Range("H126").Value = 4571.7393
$ws.Range("I126").Value = 4620.8887
$ws.Range("K126").Value = 13862.6661
$ws.Range("M126").Value = -11392.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1200.1111
$ws.Range("I81").Value = 1200.1111
$ws.Range("K81").Value = 2400.2222
$ws.Range("M81").Value = -1339.2222
$ws.Range("H84").Value = 1200.1111
$ws.Range("I84").Value = 1200.1111
$ws.Range("K84").Value = 12001.111
$ws.Range("M84").Value = -6697.111000000001
$ws.Range("I136").Value = 14109078
$ws.Range("J136").Value = 8331.666999999999
$ws.Range("K136").Value = 42327234
$ws.Range("L136").Value = 24995.001
$ws.Range("M136").Value = -42324684
$ws.Range("N136").Value = -30095.001
